$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: xyz123 / Abhishek Sisodiyha / 2232
$ws.Cells.Item(12, 1).Value = "xyz123"
$ws.Cells.Item(12, 2).Value = "Abhishek Sisodiyha"
$ws.Cells.Item(12, 3).NumberFormat = "@"
$ws.Cells.Item(12, 3).Value = "2232"
$ws.Cells.Item(12, 3).Style = "Normal"

# Row 13: vijay123 / Vijay / 1234
$ws.Cells.Item(13, 1).Value = "vijay123"
$ws.Cells.Item(13, 2).Value = "Vijay"
$ws.Cells.Item(13, 3).NumberFormat = "@"
$ws.Cells.Item(13, 3).Value = "1234"
$ws.Cells.Item(13, 3).Style = "Normal"
